$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, shifting existing rows 193:317 down to 194:318
$ws.Rows("193:193").Insert()

# The newly inserted row 193 inherits formatting from the row above/below via the
# Insert shift; now populate it with the new record's data. Columns that are not
# explicitly called out in the change (A, B, C, E, F, G, H, I, N, Q, R) keep the
# same constant values used throughout this dataset.
$ws.Range("A193").Value = 7
$ws.Range("B193").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C193").Value = "Ñuble"
$ws.Range("D193").Value = 44824
$ws.Range("E193").Value = 16
$ws.Range("F193").Value = 100112008
$ws.Range("G193").Value = "Coliflor"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 120
$ws.Range("K193").Value = 1200
$ws.Range("L193").Value = 1300
$ws.Range("M193").Value = 1250
$ws.Range("N193").Value = "$/unidad"
$ws.Range("O193").Value = "Provincia de Diguillín"
$ws.Range("P193").Value = 1250
$ws.Range("Q193").Value = 1
$ws.Range("R193").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D
$ws.Range("D193").NumberFormat = $ws.Range("D194").NumberFormat
